$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to Text format so numeric-looking strings
# (e.g. "579.93", "0.550") are stored verbatim instead of being
# auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = "65.679.41"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.482.91"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "579.93"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "160.77"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.483.70"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  +4.23%  "
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").Value = "0.126"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "4.084.04"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "0.0000195"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "28.74"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "65.623.20"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "3.479.87"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "14.29"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "390.96"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("D23").Value = "0.550"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").Value = "73.60"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "0.0000124"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "0.179"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +5.82%  "
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "23.73"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "6.48"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "7.10"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("D38").Value = "162.77"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").Value = "1.96"
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("D40").Value = "3.075.24"
$ws.Range("E40").Value = "  +5.59%  "
$ws.Range("D41").Value = "0.0770"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "27.22"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("D43").Value = "0.0320"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "4.52"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").Value = "42.75"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "0.777"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "25.91"
$ws.Range("E47").Value = "  +8.61%  "
$ws.Range("D48").Value = "1.12"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("D49").Value = "2.23"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").Value = "6.69"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").Value = "310.85"
$ws.Range("E51").Value = "  +4.32%  "
